# "revised up to V-B"
# The "Desired planning state" textbox (inside the top-level diagram group on
# slide 1) is re-worded to "Next planning state" -- split across two runs
# ("Next planning " + "state", the latter keeping the pre-existing run's
# formatting) -- and the shape is narrowed to fit the new (shorter) caption.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The shape lives inside "Group 22" (top-level shape 1); "TextBox 49"
# (shape id 50) is the 17th item of that group.
$grp = $s.Shapes.Item(1)
$shape = $grp.GroupItems.Item(17)

$tr = $shape.TextFrame.TextRange

# Replace the text, then re-split the leading "Next planning " portion back
# out into its own run so both pieces keep the shape's original Times New
# Roman run formatting (latin/cs typeface + panose/pitchFamily/charset).
$tr.Text = "Next planning state"
$lead = $tr.Characters(1, 14)
$lead.Text = "Next planning "

# Narrow the textbox from 2249334 EMU to 1980029 EMU (height is untouched).
# Shape.Width/.Left/.Top/.Height are expressed in points (1 pt = 12700 EMU);
# nudge by +0.5 EMU before converting so the EMU value round-trips exactly
# instead of drifting down by 1 EMU through the point conversion.
$targetCx = 1980029
$shape.Width = ($targetCx + 0.5) / 914400 * 72
